# Revert "1) Added Javadoc"
# This reverts the change that set C7 (password for uid 5) to "NEWPASSWORD"
# and added the extra "password1234"/"NEWPASSWORD" shared strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C7").Value = "Password"
